$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch cell outside the used range to build the replacement text as a
# genuine text value (via a string-formula), then copy/paste-special just the
# value into each BF2:BF31 cell. This avoids Excel's automatic "looks like a
# date" conversion that a direct .Value assignment of "2014-06-19" would
# trigger (turning it into a date serial number instead of literal text).
$ws.Range("BZ1").Formula = "=""2014-06-19"""

for ($i = 2; $i -le 31; $i++) {
    $ws.Range("BZ1").Copy()
    $ws.Range("BF$i").PasteSpecial(-4163)
}

# Clean up the scratch cell so it doesn't linger in the saved workbook.
$ws.Range("BZ1").ClearContents()
